$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.423.83'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '2.100.61'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = "'334.86"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +1.57%  '
$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").Value = "'0.5223"
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("D8").Value = "'0.4554"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  +3.66%  '
$ws.Range("D9").Value = "'55.91"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +12.04%  '
$ws.Range("D10").Value = "'0.08935"
$ws.Range("D10").Style = 'Normal'
$ws.Range("D11").Value = "'1.177"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("D12").Value = "'24.18"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -2.58%  '
$ws.Range("D13").Value = '2.092.77'
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("D14").Value = "'6.833"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +1.31%  '
$ws.Range("D15").Value = "'8.042"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +3.67%  '
$ws.Range("D16").Value = "'97.30"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("D17").Value = "'0.00001162"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +2.91%  '
$ws.Range("D18").Value = "'1.004"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Value = "'0.06636"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").Value = "'19.18"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -0.76%  '
$ws.Range("D21").Value = "'1.003"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").Value = "'6.298"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("D23").Value = '30.470.43'
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("D26").Value = '2.340.17'
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("D28").Value = "'163.22"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("D29").Value = "'2.517"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -4.24%  '
$ws.Range("D30").Value = "'133.38"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("D31").Value = "'1.212"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("D32").Value = "'0.1070"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("D33").Value = "'1.659"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("D34").Value = "'6.349"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +1.87%  '
$ws.Range("D35").Value = "'3.952"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +1.43%  '
$ws.Range("D36").Value = "'5.957"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +7.99%  '
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("D39").Value = "'0.06850"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +1.42%  '
$ws.Range("D40").Value = "'0.2332"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +2.39%  '
$ws.Range("D41").Value = "'12.63"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -0.72%  '
$ws.Range("D42").Value = "'0.6879"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -0.57%  '
$ws.Range("D43").Value = "'1.246"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -2.31%  '
$ws.Range("D44").Value = "'2.330"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +4.56%  '
$ws.Range("D45").Value = "'0.6412"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("E46").Value = '  -0.31%  '
$ws.Range("D47").Value = "'3.662"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("D48").Value = "'1.248"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("D49").Value = "'0.00000000346"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +15.77%  '
$ws.Range("D50").Value = "'1.204"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -1.90%  '
$ws.Range("D51").Value = "'83.04"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +0.29%  '
